$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D updates (rows 2-18)
$ws.Range("D2").Value = 105.85
$ws.Range("D3").Value = 105.61
$ws.Range("D4").Value = 105.37
$ws.Range("D5").Value = 105.13
$ws.Range("D6").Value = 104.89
$ws.Range("D7").Value = 183.3
$ws.Range("D9").Value = 985.61
$ws.Range("D10").Value = 285.65
$ws.Range("D11").Value = 142.27
$ws.Range("D12").Value = 110.43
$ws.Range("D13").Value = 111.17
$ws.Range("D14").Value = 110.15
$ws.Range("D15").Value = 110.49
$ws.Range("D16").Value = 110.96
$ws.Range("D17").Value = -191.33
$ws.Range("D18").Value = -768.21

# Column F updates (rows 2-26)
$ws.Range("F2").Value = 56.35
$ws.Range("F3").Value = 56.12
$ws.Range("F4").Value = 55.9
$ws.Range("F5").Value = 55.67
$ws.Range("F6").Value = 55.45
$ws.Range("F7").Value = 133.13
$ws.Range("F9").Value = 885.63
$ws.Range("F10").Value = 233.98
$ws.Range("F11").Value = 91.17
$ws.Range("F12").Value = 60.26
$ws.Range("F13").Value = 60.19
$ws.Range("F14").Value = 60.13
$ws.Range("F15").Value = 60.06
$ws.Range("F16").Value = 59.99
$ws.Range("F17").Value = 59.93
$ws.Range("F18").Value = 59.86
$ws.Range("F19").Value = 90.90000000000001
$ws.Range("F20").Value = 501.5
$ws.Range("F21").Value = 912.05
$ws.Range("F22").Value = 939.21
$ws.Range("F23").Value = 760.1900000000001
$ws.Range("F24").Value = 581.16
$ws.Range("F25").Value = 402.12
$ws.Range("F26").Value = 223.07
